$d = $word.ActiveDocument

$replacements = @(
    @("947÷9=105, 2", "118÷9=13, 1"),
    @("432÷9=48, 0", "479÷7=68, 3"),
    @("170÷2=85, 0", "147÷8=18, 3"),
    @("418÷5=83, 3", "490÷3=163, 1"),
    @("926÷9=102, 8", "953÷4=238, 1"),
    @("600÷3=200, 0", "271÷9=30, 1"),
    @("388÷7=55, 3", "153÷7=21, 6"),
    @("920÷7=131, 3", "200÷7=28, 4"),
    @("462÷8=57, 6", "402÷5=80, 2"),
    @("118÷4=29, 2", "951÷8=118, 7"),
    @("471÷6=78, 3", "893÷4=223, 1"),
    @("950÷8=118, 6", "153÷6=25, 3"),
    @("366÷9=40, 6", "170÷3=56, 2"),
    @("649÷9=72, 1", "757÷3=252, 1"),
    @("444÷8=55, 4", "493÷5=98, 3"),
    @("751÷4=187, 3", "650÷8=81, 2"),
    @("273÷9=30, 3", "642÷5=128, 2"),
    @("204÷7=29, 1", "187÷8=23, 3"),
    @("211÷2=105, 1", "103÷4=25, 3"),
    @("570÷7=81, 3", "783÷6=130, 3"),
    @("928÷3=309, 1", "467÷4=116, 3"),
    @("825÷7=117, 6", "403÷3=134, 1"),
    @("329÷5=65, 4", "300÷2=150, 0"),
    @("394÷2=197, 0", "226÷9=25, 1"),
    @("793÷6=132, 1", "445÷2=222, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
